$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

$ws.Range("D2").Value = 28022
$ws.Range("E2").Value = 1524
$ws.Range("F2").Value = 1524
$ws.Range("G2").Value = 719
$ws.Range("H2").Value = 433
$ws.Range("I2").Value = 695
$ws.Range("J2").Value = -262
$ws.Range("K2").Value = 36868
$ws.Range("L2").Value = 23857
$ws.Range("M2").Value = 13010
$ws.Range("N2").Value = 12173
$ws.Range("O2").Value = 837
$ws.Range("P2").Value = 1817
$ws.Range("Q2").Value = 1895
$ws.Range("R2").Value = -844
$ws.Range("S2").Value = -768
$ws.Range("T2").Value = 1338
$ws.Range("U2").Value = 557
$ws.Range("V2").Value = 17096
$ws.Range("W2").Value = 5.44
$ws.Range("X2").Value = 1.54
$ws.Range("Y2").Value = 5.79
$ws.Range("Z2").Value = 1.19
$ws.Range("AA2").Value = 183.38
$ws.Range("AB2").Value = 572.96
$ws.Range("AC2").Value = 1916
$ws.Range("AD2").Value = 13.89
$ws.Range("AE2").Value = 33759
$ws.Range("AF2").Value = 0.79
$ws.Range("AG2").Value = 550
$ws.Range("AH2").Value = 2.07
$ws.Range("AI2").Value = 28.54
$ws.Range("AJ2").Value = 36347513
$ws.Range("D3").Value = 24657
$ws.Range("E3").Value = 2179
$ws.Range("F3").Value = 2181
$ws.Range("G3").Value = 4155
$ws.Range("H3").Value = 2456
$ws.Range("I3").Value = 2722
$ws.Range("J3").Value = -266
$ws.Range("K3").Value = 36882
$ws.Range("L3").Value = 21319
$ws.Range("M3").Value = 15563
$ws.Range("N3").Value = 14439
$ws.Range("O3").Value = 1124
$ws.Range("P3").Value = 1840
$ws.Range("Q3").Value = 1085
$ws.Range("R3").Value = -962
$ws.Range("S3").Value = -264
$ws.Range("T3").Value = 681
$ws.Range("U3").Value = 404
$ws.Range("V3").Value = 15424
$ws.Range("W3").Value = 8.84
$ws.Range("X3").Value = 9.96
$ws.Range("Y3").Value = 20.46
$ws.Range("Z3").Value = 6.66
$ws.Range("AA3").Value = 136.98
$ws.Range("AB3").Value = 682.21
$ws.Range("AC3").Value = 7426
$ws.Range("AD3").Value = 4.55
$ws.Range("AE3").Value = 39550
$ws.Range("AF3").Value = 0.85
$ws.Range("AG3").Value = 750
$ws.Range("AH3").Value = 2.22
$ws.Range("AI3").Value = 10.06
$ws.Range("AJ3").Value = 36797148
$ws.Range("D4").Value = 23594
$ws.Range("E4").Value = 1493
$ws.Range("F4").Value = 1493
$ws.Range("G4").Value = 1196
$ws.Range("H4").Value = 293
$ws.Range("I4").Value = 428
$ws.Range("J4").Value = -136
$ws.Range("K4").Value = 35288
$ws.Range("L4").Value = 20049
$ws.Range("M4").Value = 15239
$ws.Range("N4").Value = 14059
$ws.Range("O4").Value = 1180
$ws.Range("P4").Value = 1877
$ws.Range("Q4").Value = 2491
$ws.Range("R4").Value = -39
$ws.Range("S4").Value = -1818
$ws.Range("T4").Value = 803
$ws.Range("U4").Value = 1688
$ws.Range("V4").Value = 13485
$ws.Range("W4").Value = 6.33
$ws.Range("X4").Value = 1.24
$ws.Range("Y4").Value = 3.01
$ws.Range("Z4").Value = 0.81
$ws.Range("AA4").Value = 131.56
$ws.Range("AB4").Value = 685.31
$ws.Range("AC4").Value = 1150
$ws.Range("AD4").Value = 28.7
$ws.Range("AE4").Value = 39069
$ws.Range("AF4").Value = 0.84
$ws.Range("AG4").Value = 750
$ws.Range("AH4").Value = 2.27
$ws.Range("AI4").Value = 63.09
$ws.Range("AJ4").Value = 37534555
$ws.Range("D5").Value = 26535
$ws.Range("E5").Value = 1757
$ws.Range("F5").Value = 1757
$ws.Range("G5").Value = 1807
$ws.Range("H5").Value = 1363
$ws.Range("I5").Value = 1100
$ws.Range("J5").Value = 263
$ws.Range("K5").Value = 36673
$ws.Range("L5").Value = 20737
$ws.Range("M5").Value = 15936
$ws.Range("N5").Value = 14500
$ws.Range("O5").Value = 1436
$ws.Range("P5").Value = 1877
$ws.Range("Q5").Value = 2093
$ws.Range("R5").Value = -1757
$ws.Range("S5").Value = 380
$ws.Range("T5").Value = 1203
$ws.Range("U5").Value = 891
$ws.Range("V5").Value = 14146
$ws.Range("W5").Value = 6.62
$ws.Range("X5").Value = 5.14
$ws.Range("Y5").Value = 7.7
$ws.Range("Z5").Value = 3.79
$ws.Range("AA5").Value = 130.13
$ws.Range("AB5").Value = 734.66
$ws.Range("AC5").Value = 2931
$ws.Range("AD5").Value = 16.03
$ws.Range("AE5").Value = 40937
$ws.Range("AF5").Value = 1.15
$ws.Range("AG5").Value = 900
$ws.Range("AH5").Value = 1.91
$ws.Range("AI5").Value = 28.97
$ws.Range("AJ5").Value = 37534555
$ws.Range("D6").Value = 27678
$ws.Range("E6").Value = 2011
$ws.Range("F6").Value = 2011
$ws.Range("G6").Value = 1821
$ws.Range("H6").Value = 1410
$ws.Range("I6").Value = 1206
$ws.Range("K6").Value = 38331
$ws.Range("L6").Value = 21383
$ws.Range("M6").Value = 16948
$ws.Range("N6").Value = 15327
$ws.Range("P6").Value = 1877
$ws.Range("Q6").Value = 2316
$ws.Range("R6").Value = -2128
$ws.Range("S6").Value = -348
$ws.Range("T6").Value = 1847
$ws.Range("U6").Value = 469
$ws.Range("V6").Value = 14573
$ws.Range("W6").Value = 7.27
$ws.Range("X6").Value = 5.1
$ws.Range("Y6").Value = 8.09
$ws.Range("Z6").Value = 3.76
$ws.Range("AA6").Value = 126.17
$ws.Range("AB6").Value = 780.98
$ws.Range("AC6").Value = 3213
$ws.Range("AD6").Value = 11.14
$ws.Range("AE6").Value = 43272
$ws.Range("AF6").Value = 0.83
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 2.79
$ws.Range("AI6").Value = 29.37
$ws.Range("AJ6").Value = 37534555
$ws.Range("D7").Value = 25735
$ws.Range("E7").Value = 1631
$ws.Range("G7").Value = 1065
$ws.Range("H7").Value = 806
$ws.Range("I7").Value = 740
$ws.Range("K7").Value = 41312
$ws.Range("L7").Value = 23909
$ws.Range("M7").Value = 17402
$ws.Range("N7").Value = 15696
$ws.Range("P7").Value = 1879
$ws.Range("Q7").Value = 2286
$ws.Range("R7").Value = -2510
$ws.Range("S7").Value = 1122
$ws.Range("T7").Value = 2444
$ws.Range("U7").Value = -74
$ws.Range("W7").Value = 6.34
$ws.Range("X7").Value = 3.13
$ws.Range("Y7").Value = 4.77
$ws.Range("Z7").Value = 2.02
$ws.Range("AA7").Value = 137.39
$ws.Range("AC7").Value = 1972
$ws.Range("AD7").Value = 26.42
$ws.Range("AE7").Value = 44313
$ws.Range("AF7").Value = 1.18
$ws.Range("AG7").Value = 915
$ws.Range("AH7").Value = 1.76
$ws.Range("AI7").Value = 46.42
$ws.Range("D8").Value = 31525
$ws.Range("E8").Value = 2817
$ws.Range("G8").Value = 3130
$ws.Range("H8").Value = 2460
$ws.Range("I8").Value = 1874
$ws.Range("K8").Value = 46101
$ws.Range("L8").Value = 26149
$ws.Range("M8").Value = 19950
$ws.Range("N8").Value = 17088
$ws.Range("P8").Value = 1879
$ws.Range("Q8").Value = 3302
$ws.Range("R8").Value = -4129
$ws.Range("S8").Value = 2009
$ws.Range("T8").Value = 4057
$ws.Range("U8").Value = -998
$ws.Range("W8").Value = 8.94
$ws.Range("X8").Value = 7.8
$ws.Range("Y8").Value = 11.43
$ws.Range("Z8").Value = 5.63
$ws.Range("AA8").Value = 131.07
$ws.Range("AC8").Value = 4992
$ws.Range("AD8").Value = 10.44
$ws.Range("AE8").Value = 48243
$ws.Range("AF8").Value = 1.08
$ws.Range("AG8").Value = 1046
$ws.Range("AH8").Value = 2.01
$ws.Range("AI8").Value = 20.96
$ws.Range("D9").Value = 34569
$ws.Range("E9").Value = 3387
$ws.Range("G9").Value = 2871
$ws.Range("H9").Value = 2211
$ws.Range("I9").Value = 1659
$ws.Range("K9").Value = 48533
$ws.Range("L9").Value = 26996
$ws.Range("M9").Value = 21534
$ws.Range("N9").Value = 18281
$ws.Range("P9").Value = 1879
$ws.Range("Q9").Value = 3698
$ws.Range("R9").Value = -2590
$ws.Range("S9").Value = -83
$ws.Range("T9").Value = 2582
$ws.Range("U9").Value = 1510
$ws.Range("W9").Value = 9.8
$ws.Range("X9").Value = 6.4
$ws.Range("Y9").Value = 9.38
$ws.Range("Z9").Value = 4.67
$ws.Range("AA9").Value = 125.36
$ws.Range("AC9").Value = 4421
$ws.Range("AD9").Value = 11.78
$ws.Range("AE9").Value = 51610
$ws.Range("AF9").Value = 1.01
$ws.Range("AG9").Value = 1158
$ws.Range("AH9").Value = 2.22
$ws.Range("AI9").Value = 26.19
